# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Value) {
    # Force text storage so numeric-looking strings (e.g. "534.19") are not
    # auto-converted to numbers, matching the original inline-string cell content.
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "58.544.49"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.613.76"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "534.19"
$ws.Range("E5").Value = "  -0.82%  "
Set-TextValue "D6" "142.59"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "2.619.77"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "3.078.61"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "58.497.85"
$ws.Range("E15").Value = "  -1.74%  "
Set-TextValue "D16" "20.77"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "2.611.41"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  -1.43%  "
Set-TextValue "D19" "4.42"
$ws.Range("E19").Value = "  +1.06%  "
Set-TextValue "D20" "334.50"
$ws.Range("E20").Value = "  -1.90%  "
Set-TextValue "D21" "10.13"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -2.10%  "
Set-TextValue "D23" "0.998"
$ws.Range("E23").Value = "  -0.14%  "
Set-TextValue "D24" "66.69"
$ws.Range("E24").Value = "  -0.96%  "
Set-TextValue "D25" "0.420"
$ws.Range("E25").Value = "  +2.65%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -2.71%  "
Set-TextValue "D28" "7.09"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0734"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -1.74%  "
Set-TextValue "D32" "5.96"
$ws.Range("E32").Value = "  +2.28%  "
Set-TextValue "D33" "153.61"
$ws.Range("E33").Value = "  +1.88%  "
Set-TextValue "D34" "18.93"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -1.41%  "
Set-TextValue "D38" "0.817"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("E39").Value = "  -2.60%  "
Set-TextValue "D40" "3.58"
$ws.Range("E40").Value = "  +1.26%  "
Set-TextValue "D41" "284.27"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("E42").Value = "  +0.09%  "
Set-TextValue "D43" "0.596"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -0.17%  "
Set-TextValue "D45" "0.0942"
$ws.Range("E45").Value = "  -0.66%  "
Set-TextValue "D46" "19.04"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "1.940.81"
$ws.Range("E49").Value = "  -0.04%  "
Set-TextValue "D50" "4.45"
$ws.Range("E50").Value = "  -1.05%  "
Set-TextValue "D51" "17.85"
$ws.Range("E51").Value = "  -3.35%  "
